# Brief borrow book.docx — apply the commit's edits via Word COM-interop.
$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Text corrections / rewrites (Find & Replace also collapses any
#    split runs covering the matched text into a single clean run).
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "3.2 Borrow Confirmation", $true, $false, $false, $false, $false,
    $true, 1, $false, "3.2 Borrow Confirmation", 2) | Out-Null

$d.Content.Find.Execute(
    "The system will accept searching details from the user. After that filter the all the book in the library after that will display the searching result for user.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The system will accept searching details from the user. Then filter the all the book in the library that meet the search requirement and display the result to the user.",
    2) | Out-Null

$d.Content.Find.Execute(
    "After user searching process the system will allow user select the book want to borrow. After that, system will required user to confirm the book they want to borrow.  ",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "After the searching process, the user is allowed to select one book to borrow. After that, system will require the user to confirm the book they want to borrow.",
    2) | Out-Null

$d.Content.Find.Execute(
    "After the confirmation process the detail of the user will record on the book reservation file for future needs. After that, system will display the book status for the user and inform user when to collect the book.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "After the confirmation process, the detail of the user will be recorded on the book reservation file for future use. After that, system will display the book status for the user and inform the user when to collect the book at the library.",
    2) | Out-Null

$d.Content.Find.Execute(
    "User come to library to collect the book.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "When the user come to library to collect the book, the system will change the status of the book from reserved to collected. If the user did not collect the book at the specific day, the system will automatically change the status of book to available and assume that the user has returned the book.",
    2) | Out-Null

$d.Content.Find.Execute(
    "If user not enough time to read this book the system also allow user make an extension request. The extension request process required user key in the serial book on the book and user ID. After that, system will generate a new expired date of the book for the user.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "If user not enough time to read this book, the system also allow user make an extension request. The extension request process requires the user to key in the serial number on the book and their user ID. After that, system will generate and display a new expired date of the book for the user. User are only allowed to use this feature not more that 2 times for each book they borrowed to prevent the user from abusing this feature.",
    2) | Out-Null

$d.Content.Find.Execute(
    "After user collect the book or renew the book the system will update the book status in the book file.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "After user collected the book or renewed the book, the system will update the book status in the book file.",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. Uniform paragraph line spacing (w:spacing w:line="360" w:lineRule="auto")
#    on every paragraph.
# ---------------------------------------------------------------------
$d.Paragraphs.LineSpacingRule = 5
$d.Paragraphs.LineSpacing = 18

# ---------------------------------------------------------------------
# 3. Title / "3.1 Search Book" heading mark: sz 72 (36pt) -> sz 24 (12pt).
# ---------------------------------------------------------------------
$d.Paragraphs(1).Range.Font.Size = 12
$d.Paragraphs(1).Range.Font.SizeBi = 12
$d.Paragraphs(2).Range.Font.Size = 12
$d.Paragraphs(2).Range.Font.SizeBi = 12

# ---------------------------------------------------------------------
# 4. Merge each "3.x heading" paragraph with the body paragraph that
#    follows it: append a manual line break (<w:br/>) to the end of the
#    heading paragraph, then delete the paragraph mark that used to
#    separate the two paragraphs so their runs become one paragraph.
# ---------------------------------------------------------------------
function Merge-WithNextParagraph($headingIndex) {
    $heading = $d.Paragraphs($headingIndex)
    $endPt = $d.Range($heading.Range.End - 1, $heading.Range.End - 1)
    $endPt.InsertBreak(6) | Out-Null
    $heading2 = $d.Paragraphs($headingIndex)
    $markRange = $d.Range($heading2.Range.End - 1, $heading2.Range.End)
    $markRange.Delete() | Out-Null
}

# Heading paragraph indices (1-based) before any merging: 2,4,6,8,10,12.
# Walk from the last pair to the first so earlier indices stay valid.
Merge-WithNextParagraph 12
Merge-WithNextParagraph 10
Merge-WithNextParagraph 8
Merge-WithNextParagraph 6
Merge-WithNextParagraph 4
Merge-WithNextParagraph 2

Write-Output "merges done"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Output ($i.ToString() + ": [" + $d.Paragraphs($i).Range.Text + "]")
}
